# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) / "Valor Mora" (column F) table occupying
# rows 16-53 of Hoja1 gets re-sorted from descending period order
# (1908 .. 1607) to ascending period order (1607 .. 1908), with each
# "Valor Mora" amount staying attached to its own period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{Row=16; Period="1607"; Valor=27578},
    @{Row=17; Period="1608"; Valor=27578},
    @{Row=18; Period="1609"; Valor=27578},
    @{Row=19; Period="1610"; Valor=27578},
    @{Row=20; Period="1611"; Valor=27578},
    @{Row=21; Period="1612"; Valor=27578},
    @{Row=22; Period="1701"; Valor=27578},
    @{Row=23; Period="1702"; Valor=27578},
    @{Row=24; Period="1703"; Valor=27578},
    @{Row=25; Period="1704"; Valor=27578},
    @{Row=26; Period="1705"; Valor=27578},
    @{Row=27; Period="1706"; Valor=27578},
    @{Row=28; Period="1707"; Valor=27578},
    @{Row=29; Period="1708"; Valor=27578},
    @{Row=30; Period="1709"; Valor=27578},
    @{Row=31; Period="1710"; Valor=27578},
    @{Row=32; Period="1711"; Valor=27578},
    @{Row=33; Period="1712"; Valor=27578},
    @{Row=34; Period="1801"; Valor=27578},
    @{Row=35; Period="1802"; Valor=27578},
    @{Row=36; Period="1803"; Valor=27578},
    @{Row=37; Period="1804"; Valor=27578},
    @{Row=38; Period="1805"; Valor=27578},
    @{Row=39; Period="1806"; Valor=27578},
    @{Row=40; Period="1807"; Valor=27578},
    @{Row=41; Period="1808"; Valor=27578},
    @{Row=42; Period="1809"; Valor=31249},
    @{Row=43; Period="1810"; Valor=31249},
    @{Row=44; Period="1811"; Valor=31249},
    @{Row=45; Period="1812"; Valor=31249},
    @{Row=46; Period="1901"; Valor=31249},
    @{Row=47; Period="1902"; Valor=31249},
    @{Row=48; Period="1903"; Valor=31249},
    @{Row=49; Period="1904"; Valor=31249},
    @{Row=50; Period="1905"; Valor=31249},
    @{Row=51; Period="1906"; Valor=31249},
    @{Row=52; Period="1907"; Valor=31249},
    @{Row=53; Period="1908"; Valor=30208}
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 5).Value = $item.Period
    $ws.Cells.Item($item.Row, 6).Value = $item.Valor
}
